# Add login credentials for "Script 2" to the Login sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Duplicate the existing credential block (header row + data row + blank/footer
# row) a few rows below the current one, keeping the same formatting.
$ws.Range("A1:G3").Copy($ws.Range("A11")) | Out-Null

# New section title above the duplicated block, highlighted in yellow.
$ws.Range("A10").Value = "Script 2"
$ws.Range("A10").Interior.Color = 65535

# The new block uses its own username for this script.
$ws.Range("A12").Value = "AGSAutoT03"

# Leave the selection on the newly added username cell, like the author did.
$ws.Range("A12").Select() | Out-Null
